$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("E2").Value = 3
$ws.Range("G2").Value = 1.807599666666667
$ws.Range("H2").Value = 5.422799
$ws.Range("K2").Value = 3
$ws.Range("M2").Value = 5.057757666666666
$ws.Range("N2").Value = 15.173273
$ws.Range("O2").Value = 0.173378811020062
$ws.Range("P2").Value = 0.173378811020062
$ws.Range("Q2").Value = 9.142401072347445
$ws.Range("R2").Value = 82.281609651127
$ws.Range("S2").Value = 0.173378811020062
$ws.Range("T2").Value = 0.173378811020062

# Row 3
$ws.Range("E3").Value = 3
$ws.Range("G3").Value = 1.807599666666667
$ws.Range("H3").Value = 5.422799
$ws.Range("K3").Value = 3
$ws.Range("M3").Value = 14.247411
$ws.Range("N3").Value = 42.742233
$ws.Range("O3").Value = 0.4883980890531961
$ws.Range("P3").Value = 0.4883980890531961
$ws.Range("Q3").Value = 25.753615374463
$ws.Range("R3").Value = 231.782538370167
$ws.Range("S3").Value = 0.4883980890531961
$ws.Range("T3").Value = 0.4883980890531961

# Row 4
$ws.Range("E4").Value = 3
$ws.Range("G4").Value = 1.807599666666667
$ws.Range("H4").Value = 5.422799
$ws.Range("K4").Value = 3
$ws.Range("M4").Value = 9.866548666666667
$ws.Range("N4").Value = 29.599646
$ws.Range("O4").Value = 0.3382230999267418
$ws.Range("P4").Value = 0.3382230999267418
$ws.Range("Q4").Value = 17.83477008101711
$ws.Range("R4").Value = 160.512930729154
$ws.Range("S4").Value = 0.3382230999267418
$ws.Range("T4").Value = 0.3382230999267418
